$d = $word.ActiveDocument

function Get-ParaIndexByText($doc, $text) {
    $paras = $doc.Paragraphs
    $count = $paras.Count
    for ($i = 1; $i -le $count; $i++) {
        $t = $paras.Item($i).Range.Text.TrimEnd([char]13)
        if ($t -eq $text) {
            return $i
        }
    }
    return -1
}

function Remove-ParaRange($doc, $startText, $endText) {
    $paras = $doc.Paragraphs
    $startIdx = Get-ParaIndexByText $doc $startText
    $endIdx = Get-ParaIndexByText $doc $endText
    $startRange = $paras.Item($startIdx).Range
    $endRange = $paras.Item($endIdx).Range
    $delRange = $doc.Range($startRange.Start, $endRange.End)
    $delRange.Delete()
}

# 1. Remove the extra cafe entries between "C3" and "Food Type : " block
#    (C4, Retro, Coffee Lounge, NSTP Café, NBS Café, S3H Café, Khaapa SEECS, SCME Cafe)
Remove-ParaRange $d "C4" "SCME Cafe"

# 2. Remove the extra food-type entries between "Cold" and the trailing "Sweet"
#    (Snacks, Sweet, Spicy)
Remove-ParaRange $d "Snacks" "Spicy"

# 3. Split the "Cold" run into "Col" + "d" (same formatting), matching the
#    author's final edit where the last letter became its own run.
$paras = $d.Paragraphs
$coldIdx = Get-ParaIndexByText $d "Cold"
$coldRange = $paras.Item($coldIdx).Range
$lastCharRange = $d.Range($coldRange.End - 2, $coldRange.End - 1)
$lastCharRange.Font.Bold = 1
$lastCharRange.Font.Bold = 0
